$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the small brace shape that decorated the "Laufend" row (row 8),
# which is being removed below.
$ws.Shapes.Item("Geschweifte Klammer rechts 3").Delete()

# Remove row 8 ("Laufend" / "Wenn Distanz > vorherigeDistanz + X") entirely.
$ws.Rows("8").Delete()

# After the shift above, the old "~ 3 sec. / laufend / ... / Default" row
# is now row 12; remove it too so "bis Ende / laufend-stehend / ..." moves
# up to become the new row 12.
$ws.Rows("12").Delete()

# Rename the remaining "Default" algorithm status to "Standard".
$ws.Range("G11").Value = "Standard"

# Update the view: zoom to 130% and select the Ablauf table range.
$excel.ActiveWindow.Zoom = 130
$ws.Range("C9:G12").Select()
